$wb = $excel.ActiveWorkbook
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$newSheet.Name = "Sheet5"
$newSheet.Range("A1").Value = "Test"
$newSheet.Range("D1").Formula = "=_xlfn.UNIQUE(_xlfn._xlws.SORT(_xlfn._xlws.FILTER(B:B,B:B<>"""")))"
Write-Host $newSheet.Range("D1").Value
